$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.883
$ws.Range("B8").Value = 6.467000000000001
$ws.Range("B10").Value = 6.484
$ws.Range("B12").Value = 5.553
$ws.Range("D12").Value = -6.985000000000001
$ws.Range("D15").Value = -8.312000000000001
$ws.Range("D17").Value = -8.126999999999999
$ws.Range("B18").Value = 5.137
$ws.Range("D26").Value = -7.194999999999999
$ws.Range("D27").Value = -7.923999999999999
$ws.Range("D28").Value = -8.002000000000001
$ws.Range("B37").Value = 8.882000000000001
$ws.Range("D37").Value = -7.900999999999999
$ws.Range("D47").Value = -7.568
$ws.Range("B55").Value = 4.722
$ws.Range("D65").Value = -7.659000000000001
$ws.Range("B68").Value = 5.084999999999999
$ws.Range("D73").Value = -8.023
$ws.Range("B77").Value = 5.766
$ws.Range("B78").Value = 7.507000000000001
$ws.Range("B81").Value = 6.439
$ws.Range("B82").Value = 5.457
$ws.Range("D84").Value = -8.235000000000001
$ws.Range("D85").Value = -8.574
$ws.Range("D93").Value = -7.007000000000001
$ws.Range("D95").Value = -7.558
$ws.Range("D98").Value = -7.231
$ws.Range("D99").Value = -8.269000000000002
$ws.Range("D101").Value = -7.834000000000001
